# Updated symbol list on Sat Dec 24 23:44:11 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds decimal-looking values that must stay TEXT (the sheet
# stores them as inlineStr, not numbers). Prefixing with an apostrophe
# forces Excel to keep them as literal text instead of re-parsing /
# re-formatting them as floating point numbers (which would lose
# trailing zeros, e.g. "244.70" -> 244.7).

# Simple price (column D) updates
$ws.Range("D2").Value  = "'244.70"
$ws.Range("D3").Value  = "'21.83"
$ws.Range("D4").Value  = "'5.399"
$ws.Range("D5").Value  = "'0.06045"
$ws.Range("D7").Value  = "'0.8141"
$ws.Range("D8").Value  = "'0.9239"
$ws.Range("D9").Value  = "'0.1439"
$ws.Range("D10").Value = "'0.07477"
$ws.Range("D11").Value = "'0.03380"
$ws.Range("D12").Value = "'0.03050"

# Rows 13 & 14 swap (BitMartToken <-> MCDex), with new prices
$ws.Range("B13").Value = "MCDex"
$ws.Range("C13").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D13").Value = "'4.006"
$ws.Range("E13").Value = "12MCDexMCB"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09408"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("D15").Value = "'0.001590"
$ws.Range("D16").Value = "'0.04825"

$ws.Range("D17").Value = "'0.01120"
$ws.Range("E17").Value = "16OneONEBestin24h"

$ws.Range("D18").Value = "'0.005324"

$ws.Range("D20").Value = "'0.0009899"
$ws.Range("D21").Value = "'0.00008805"
$ws.Range("D22").Value = "'3.653"
$ws.Range("D23").Value = "'6.429"
$ws.Range("D24").Value = "'2.145"

$ws.Range("D26").Value = "'0.1340"
$ws.Range("D27").Value = "'0.0002901"

$ws.Range("D40").Value = "'0.03987"

# Rows 41, 42, 43 rotate (BKEXToken -> row42, CEJI -> row43, KickToken -> row41)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006413"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1075"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002901"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.006387"
$ws.Range("D45").Value = "'0.00005239"
$ws.Range("D46").Value = "'0.00000000750"

$ws.Range("D47").Value = "'1.100"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").Value = "'0.002527"
$ws.Range("D50").Value = "'0.01010"
